$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve text formatting for numeric-looking price strings so Excel
# does not silently coerce them into Number cells.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated coin data.
$ws.Range("D2").Value = "63.792.15"
$ws.Range("E2").Value = "  -0.24%  "
$ws.Range("D3").Value = "3.077.40"
$ws.Range("E3").Value = "  +0.10%  "
$ws.Range("E4").Value = "  -0.43%  "
$ws.Range("D5").Value = "592.49"
$ws.Range("E5").Value = "  +0.81%  "
$ws.Range("D6").Value = "155.20"
$ws.Range("E6").Value = "  +1.50%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").Value = "0.539"
$ws.Range("E8").Value = "  +1.36%  "
$ws.Range("D9").Value = "3.076.79"
$ws.Range("E9").Value = "  +0.14%  "
$ws.Range("D10").Value = "0.157"
$ws.Range("E10").Value = "  -0.39%  "
$ws.Range("D11").Value = "5.91"
$ws.Range("E11").Value = "  +0.22%  "
$ws.Range("D12").Value = "0.453"
$ws.Range("E12").Value = "  -0.89%  "
$ws.Range("B13").Value = "ShibaInu"
$ws.Range("C13").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D13").Value = "0.0000239"
$ws.Range("E13").Value = "  -1.69%  "
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").Value = "36.89"
$ws.Range("E14").Value = "  -1.70%  "
$ws.Range("E15").Value = "  +1.59%  "
$ws.Range("D16").Value = "3.590.12"
$ws.Range("E16").Value = "  +0.08%  "
$ws.Range("D17").Value = "7.21"
$ws.Range("E17").Value = "  +0.80%  "
$ws.Range("D18").Value = "63.764.15"
$ws.Range("E18").Value = "  +0.12%  "
$ws.Range("D19").Value = "3.081.04"
$ws.Range("E19").Value = "  +0.26%  "
$ws.Range("D20").Value = "485.00"
$ws.Range("E20").Value = "  +3.57%  "
$ws.Range("D21").Value = "14.59"
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").Value = "0.711"
$ws.Range("E22").Value = "  -2.10%  "
$ws.Range("D23").Value = "7.62"
$ws.Range("E23").Value = "  +1.37%  "
$ws.Range("E24").Value = "  +3.23%  "
$ws.Range("D25").Value = "82.00"
$ws.Range("E25").Value = "  +1.11%  "
$ws.Range("D26").Value = "12.91"
$ws.Range("E26").Value = "  -2.50%  "
$ws.Range("E27").Value = "  +8.14%  "
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.16%  "
$ws.Range("D29").Value = "7.52"
$ws.Range("E29").Value = "  +2.90%  "
$ws.Range("D30").Value = "2.26"
$ws.Range("E30").Value = "  +3.02%  "
$ws.Range("D31").Value = "2.70"
$ws.Range("E31").Value = "  +0.73%  "
$ws.Range("D32").Value = "1.00"
$ws.Range("E32").Value = "  -0.41%  "
$ws.Range("D33").Value = "27.33"
$ws.Range("E33").Value = "  -0.05%  "
$ws.Range("D34").Value = "0.112"
$ws.Range("E34").Value = "  -2.39%  "
$ws.Range("B35").Value = "Mantle"
$ws.Range("C35").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D35").Value = "1.07"
$ws.Range("E35").Value = "  +1.75%  "
$ws.Range("B36").Value = "PEPE"
$ws.Range("C36").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D36").Value = "0.0₃0824"
$ws.Range("E36").Value = "  -3.22%  "
$ws.Range("D37").Value = "6.08"
$ws.Range("E37").Value = "  -0.27%  "
$ws.Range("D38").Value = "2.24"
$ws.Range("E38").Value = "  -0.14%  "
$ws.Range("D39").Value = "3.25"
$ws.Range("E39").Value = "  -4.00%  "
$ws.Range("D40").Value = "9.31"
$ws.Range("E40").Value = "  +0.37%  "
$ws.Range("D41").Value = "50.76"
$ws.Range("E41").Value = "  +0.35%  "
$ws.Range("D42").Value = "440.44"
$ws.Range("E42").Value = "  -1.54%  "
$ws.Range("D43").Value = "0.291"
$ws.Range("E43").Value = "  +1.50%  "
$ws.Range("D44").Value = "0.0366"
$ws.Range("E44").Value = "  +0.38%  "
$ws.Range("E45").Value = "  +3.97%  "
$ws.Range("D46").Value = "2.839.70"
$ws.Range("E46").Value = "  +0.64%  "
$ws.Range("D47").Value = "39.50"
$ws.Range("E47").Value = "  -0.81%  "
$ws.Range("D48").Value = "131.71"
$ws.Range("E48").Value = "  +1.80%  "
$ws.Range("D49").Value = "25.60"
$ws.Range("E49").Value = "  +0.96%  "
$ws.Range("E50").Value = "  +0.01%  "
$ws.Range("D51").Value = "2.25"
$ws.Range("E51").Value = "  +0.69%  "
